$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'89.544.65"
$ws.Range("E2").Value = '  +0.37%  '

$ws.Range("D3").Value = "'3.069.49"
$ws.Range("E3").Value = '  -1.93%  '

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").Value = "'212.16"
$ws.Range("E5").Value = '  +0.15%  '

$ws.Range("D6").Value = "'610.27"
$ws.Range("E6").Value = '  -2.35%  '

$ws.Range("E7").Value = '  -7.09%  '

$ws.Range("D8").Value = "'0.907"
$ws.Range("E8").Value = '  +27.06%  '

$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = '  +0.11%  '

$ws.Range("D10").Value = "'3.070.27"
$ws.Range("E10").Value = '  -1.78%  '

$ws.Range("D11").Value = "'0.672"
$ws.Range("E11").Value = '  +22.40%  '

$ws.Range("E12").Value = '  +4.93%  '

$ws.Range("D13").Value = "'0.0000239"
$ws.Range("E13").Value = '  -4.31%  '

$ws.Range("D14").Value = "'5.39"
$ws.Range("E14").Value = '  +2.84%  '

$ws.Range("D15").Value = "'89.391.58"
$ws.Range("E15").Value = '  +0.27%  '

$ws.Range("D16").Value = "'32.24"
$ws.Range("E16").Value = '  +1.48%  '

$ws.Range("D17").Value = "'3.644.93"
$ws.Range("E17").Value = '  -1.78%  '

$ws.Range("D18").Value = "'3.089.06"
$ws.Range("E18").Value = '  -2.06%  '

$ws.Range("D19").Value = "'3.38"
$ws.Range("E19").Value = '  +2.56%  '

$ws.Range("D20").Value = "'0.0000209"
$ws.Range("E20").Value = '  -0.70%  '

$ws.Range("D21").Value = "'13.51"
$ws.Range("E21").Value = '  +3.07%  '

$ws.Range("D22").Value = "'428.75"
$ws.Range("E22").Value = '  +1.55%  '

$ws.Range("E23").Value = '  +2.58%  '

$ws.Range("E24").Value = '  -1.78%  '

$ws.Range("D25").Value = "'5.44"
$ws.Range("E25").Value = '  +5.32%  '

$ws.Range("D26").Value = "'84.61"
$ws.Range("E26").Value = '  +7.03%  '

$ws.Range("D27").Value = "'11.79"
$ws.Range("E27").Value = '  +3.48%  '

$ws.Range("D28").Value = "'3.236.53"
$ws.Range("E28").Value = '  -3.44%  '

$ws.Range("E29").Value = '  +0.06%  '

$ws.Range("E30").Value = '  +10.00%  '

$ws.Range("E31").Value = '  +4.56%  '

$ws.Range("D32").Value = "'8.25"
$ws.Range("E32").Value = '  +0.94%  '

$ws.Range("D33").Value = "'505.03"
$ws.Range("E33").Value = '  -0.12%  '

$ws.Range("D34").Value = "'3.59"
$ws.Range("E34").Value = '  -7.69%  '

$ws.Range("D35").Value = "'6.63"
$ws.Range("E35").Value = '  -1.71%  '

$ws.Range("D36").Value = "'22.79"
$ws.Range("E36").Value = '  +4.27%  '

$ws.Range("B37").Value = 'PancakeSwap'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D37").Value = "'1.79"
$ws.Range("E37").Value = '  -3.43%  '

$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").Value = "'1.24"
$ws.Range("E38").Value = '  -1.43%  '

$ws.Range("D39").Value = "'0.132"
$ws.Range("E39").Value = '  +5.25%  '

$ws.Range("D40").Value = "'22.26"
$ws.Range("E40").Value = '  +0.06%  '

$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = '  +0.43%  '

$ws.Range("E42").Value = '  +0.06%  '

$ws.Range("D43").Value = "'0.143"
$ws.Range("E43").Value = '  +15.88%  '

$ws.Range("D44").Value = "'0.368"
$ws.Range("E44").Value = '  +1.67%  '

$ws.Range("E45").Value = '  -1.83%  '

$ws.Range("D46").Value = "'147.83"
$ws.Range("E46").Value = '  +0.74%  '

$ws.Range("D47").Value = "'0.0693"
$ws.Range("E47").Value = '  +14.73%  '

$ws.Range("D48").Value = "'43.46"
$ws.Range("E48").Value = '  -0.20%  '

$ws.Range("D49").Value = "'4.09"
$ws.Range("E49").Value = '  +4.35%  '

$ws.Range("E50").Value = '  +2.97%  '

$ws.Range("D51").Value = "'156.18"
$ws.Range("E51").Value = '  -6.07%  '
